$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Shift the widths of columns C, D, E (3,4,5) one slot to the right ---
# Before: C=49  D=51  E=37
# After : C=37  D=49  E=51
$w3 = $ws.Columns.Item(3).ColumnWidth
$w4 = $ws.Columns.Item(4).ColumnWidth
$w5 = $ws.Columns.Item(5).ColumnWidth

$ws.Columns.Item(3).ColumnWidth = $w5
$ws.Columns.Item(4).ColumnWidth = $w3
$ws.Columns.Item(5).ColumnWidth = $w4

# --- Row 1 (header): insert a new value in C1, shift D1/E1 right, new value in F1 ---
$ws.Range("C1").Value = "button_testResultActions_class_2"
$ws.Range("D1").Value = "button_testResultActions_internalRoleButtonName"
$ws.Range("E1").Value = "button_testResultActions_internalRoleButtonName_1"
$ws.Range("F1").Value = "button_testResultActions_internalRoleButtonName_2"

# --- Row 2 (data): move old E2 into C2, shift C2/D2 right into D2/E2, F2 unchanged ---
$ws.Range("C2").Value = '"]:nth-child(3) [class="css-1yjo05o'
$ws.Range("D2").Value = "Failed Automations - Apply to"
$ws.Range("E2").Value = "Failed Portal - Login with"
